$wb = $excel.ActiveWorkbook

# Row 108 on ALC (hunk 0)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(108, 8).Value = 69295
$ws.Cells.Item(108, 10).Value = 69295
$ws.Cells.Item(108, 12).Value = 69295
$ws.Cells.Item(108, 14).Value = -76975

# Row 132 on ALC (hunk 1)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1506.0303
$ws.Cells.Item(132, 9).Value = 1520.1
$ws.Cells.Item(132, 11).Value = 4560.299999999999
$ws.Cells.Item(132, 13).Value = -2030.299999999999

# Row 133 on ALC (hunk 2)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(133, 8).Value = 98995.664
$ws.Cells.Item(133, 10).Value = 98995.664
$ws.Cells.Item(133, 12).Value = 98995.664
$ws.Cells.Item(133, 14).Value = -109115.664

# Row 134 on ALC (hunk 3)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(134, 8).Value = 87648.336
$ws.Cells.Item(134, 10).Value = 87648.336
$ws.Cells.Item(134, 12).Value = 87648.336
$ws.Cells.Item(134, 14).Value = -97788.336

# Row 136 on ALC (hunk 4)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(136, 8).Value = 78605
$ws.Cells.Item(136, 10).Value = 78605
$ws.Cells.Item(136, 12).Value = 78605
$ws.Cells.Item(136, 14).Value = -88805

# Row 138 on ALC (hunk 5)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2351.6667
$ws.Cells.Item(138, 9).Value = 1691.3889
$ws.Cells.Item(138, 10).Value = 2917.6191
$ws.Cells.Item(138, 11).Value = 5074.1667
$ws.Cells.Item(138, 12).Value = 8752.8573
$ws.Cells.Item(138, 13).Value = 65.83330000000024
$ws.Cells.Item(138, 14).Value = -19032.8573

# Row 139 on ALC (hunk 6)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(139, 8).Value = 70137.91
$ws.Cells.Item(139, 10).Value = 70137.91
$ws.Cells.Item(139, 12).Value = 70137.91
$ws.Cells.Item(139, 14).Value = -80417.91

# Row 140 on ALC (hunk 7)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(140, 8).Value = 87617.2
$ws.Cells.Item(140, 10).Value = 87617.2
$ws.Cells.Item(140, 12).Value = 87617.2
$ws.Cells.Item(140, 14).Value = -97977.2

# Row 141 on ALC (hunk 8)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 2396.3076
$ws.Cells.Item(141, 9).Value = 2187.2173
$ws.Cells.Item(141, 10).Value = 3999.3333
$ws.Cells.Item(141, 11).Value = 6561.651899999999
$ws.Cells.Item(141, 12).Value = 11997.9999
$ws.Cells.Item(141, 13).Value = -1381.651899999999
$ws.Cells.Item(141, 14).Value = -22357.9999

# Row 32 on ARM (hunk 9)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5530.068
$ws.Cells.Item(32, 9).Value = 2487.375
$ws.Cells.Item(32, 10).Value = 18807.273
$ws.Cells.Item(32, 11).Value = 2487.375
$ws.Cells.Item(32, 12).Value = 18807.273
$ws.Cells.Item(32, 13).Value = -2200.375
$ws.Cells.Item(32, 14).Value = -19381.273

# Row 52 on ARM (hunk 10)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(52, 8).Value = 56498
$ws.Cells.Item(52, 10).Value = 56498
$ws.Cells.Item(52, 12).Value = 56498
$ws.Cells.Item(52, 14).Value = -57134

# Row 63 on ARM (hunk 11)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 3619
$ws.Cells.Item(63, 9).Value = 2094.5
$ws.Cells.Item(63, 11).Value = 2094.5
$ws.Cells.Item(63, 13).Value = -1408.5

# Row 66 on ARM (hunk 12)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 3619
$ws.Cells.Item(66, 9).Value = 2094.5
$ws.Cells.Item(66, 11).Value = 10472.5
$ws.Cells.Item(66, 13).Value = -7040.5

# Row 97 on ARM (hunk 13)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 643.6667
$ws.Cells.Item(97, 9).Value = 666.4666999999999
$ws.Cells.Item(97, 10).Value = 586.6667
$ws.Cells.Item(97, 11).Value = 666.4666999999999
$ws.Cells.Item(97, 12).Value = 586.6667
$ws.Cells.Item(97, 13).Value = -170.4666999999999
$ws.Cells.Item(97, 14).Value = -1578.6667

# Row 117 on ARM (hunk 14)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(117, 8).Value = 66607.71000000001
$ws.Cells.Item(117, 10).Value = 66607.71000000001
$ws.Cells.Item(117, 12).Value = 66607.71000000001
$ws.Cells.Item(117, 14).Value = -75785.71000000001

# Row 118 on ARM (hunk 15)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(118, 8).Value = 45809.832
$ws.Cells.Item(118, 10).Value = 45809.832
$ws.Cells.Item(118, 12).Value = 45809.832
$ws.Cells.Item(118, 14).Value = -49123.832

# Row 132 on ARM (hunk 16)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 3437.2432
$ws.Cells.Item(132, 10).Value = 7505.8335
$ws.Cells.Item(132, 12).Value = 22517.5005
$ws.Cells.Item(132, 14).Value = -27577.5005

# Row 55 on BSM (hunk 17)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(55, 8).Value = 27940.428
$ws.Cells.Item(55, 10).Value = 27940.428
$ws.Cells.Item(55, 12).Value = 27940.428
$ws.Cells.Item(55, 14).Value = -28486.428

# Row 115 on BSM (hunk 18)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(115, 8).Value = 77664.44500000001
$ws.Cells.Item(115, 10).Value = 77664.44500000001
$ws.Cells.Item(115, 12).Value = 77664.44500000001
$ws.Cells.Item(115, 14).Value = -80798.44500000001

# Row 134 on BSM (hunk 19)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 6141.185
$ws.Cells.Item(134, 9).Value = 3611.7222
$ws.Cells.Item(134, 11).Value = 10835.1666
$ws.Cells.Item(134, 13).Value = -8300.1666

# Row 135 on BSM (hunk 20)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(135, 8).Value = 99996
$ws.Cells.Item(135, 10).Value = 99996
$ws.Cells.Item(135, 12).Value = 99996
$ws.Cells.Item(135, 14).Value = -110136

# Row 138 on BSM (hunk 21)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(138, 8).Value = 99994.164
$ws.Cells.Item(138, 10).Value = 99994.164
$ws.Cells.Item(138, 12).Value = 99994.164
$ws.Cells.Item(138, 14).Value = -110274.164

# Row 118 on CRP (hunk 22)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(118, 8).Value = 89995.57000000001
$ws.Cells.Item(118, 10).Value = 89995.57000000001
$ws.Cells.Item(118, 12).Value = 89995.57000000001
$ws.Cells.Item(118, 14).Value = -93309.57000000001

# Row 134 on CRP (hunk 23)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 8932561
$ws.Cells.Item(134, 9).Value = 11908957
$ws.Cells.Item(134, 11).Value = 35726871
$ws.Cells.Item(134, 13).Value = -35724336

# Row 138 on CRP (hunk 24)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(138, 8).Value = 99986.664
$ws.Cells.Item(138, 10).Value = 99986.664
$ws.Cells.Item(138, 12).Value = 99986.664
$ws.Cells.Item(138, 14).Value = -110266.664

# Row 107 on CUL (hunk 25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 1857.1428
$ws.Cells.Item(107, 9).Value = 1857.1428
$ws.Cells.Item(107, 11).Value = 5571.428400000001
$ws.Cells.Item(107, 13).Value = -3651.428400000001

# Row 108 on CUL (hunk 26)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(108, 8).Value = 136.28572
$ws.Cells.Item(108, 10).Value = 125
$ws.Cells.Item(108, 12).Value = 375
$ws.Cells.Item(108, 14).Value = -6135

# Row 110 on CUL (hunk 27)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(110, 8).Value = 4587.5713
$ws.Cells.Item(110, 9).Value = 3410.6
$ws.Cells.Item(110, 10).Value = 7530
$ws.Cells.Item(110, 11).Value = 10231.8
$ws.Cells.Item(110, 12).Value = 22590
$ws.Cells.Item(110, 13).Value = -6141.799999999999
$ws.Cells.Item(110, 14).Value = -30770

# Row 121 on CUL (hunk 28)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(121, 8).Value = 4376326
$ws.Cells.Item(121, 10).Value = 17500000
$ws.Cells.Item(121, 12).Value = 52500000
$ws.Cells.Item(121, 14).Value = -52502620

# Row 131 on CUL (hunk 29)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1720.625
$ws.Cells.Item(131, 9).Value = 1248.75
$ws.Cells.Item(131, 10).Value = 2192.5
$ws.Cells.Item(131, 11).Value = 3746.25
$ws.Cells.Item(131, 12).Value = 6577.5
$ws.Cells.Item(131, 13).Value = 1293.75
$ws.Cells.Item(131, 14).Value = -16657.5

# Row 10 on GSM (hunk 30)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(10, 8).Value = 276810.1
$ws.Cells.Item(10, 9).Value = 2234.3333
$ws.Cells.Item(10, 10).Value = 606301
$ws.Cells.Item(10, 11).Value = 2234.3333
$ws.Cells.Item(10, 12).Value = 606301
$ws.Cells.Item(10, 13).Value = -2065.3333
$ws.Cells.Item(10, 14).Value = -606639

# Row 11 on GSM (hunk 31)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 3032762.2
$ws.Cells.Item(11, 9).Value = 1464627
$ws.Cells.Item(11, 10).Value = 4339541.5
$ws.Cells.Item(11, 11).Value = 1464627
$ws.Cells.Item(11, 12).Value = 4339541.5
$ws.Cells.Item(11, 13).Value = -1464488
$ws.Cells.Item(11, 14).Value = -4339819.5

# Row 12 on GSM (hunk 32)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(12, 8).Value = 4271.12
$ws.Cells.Item(12, 9).Value = 4898.706
$ws.Cells.Item(12, 11).Value = 4898.706
$ws.Cells.Item(12, 13).Value = -4758.706

# Row 109 on GSM (hunk 33)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(109, 8).Value = 36679.69
$ws.Cells.Item(109, 10).Value = 36679.69
$ws.Cells.Item(109, 12).Value = 36679.69
$ws.Cells.Item(109, 14).Value = -38759.69

# Row 110 on GSM (hunk 34)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(110, 8).Value = 99896.664
$ws.Cells.Item(110, 10).Value = 99896.664
$ws.Cells.Item(110, 12).Value = 99896.664
$ws.Cells.Item(110, 14).Value = -108076.664

# Row 119 on GSM (hunk 35)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(119, 8).Value = 85830.234
$ws.Cells.Item(119, 10).Value = 85830.234
$ws.Cells.Item(119, 12).Value = 85830.234
$ws.Cells.Item(119, 14).Value = -95506.234

# Row 135 on GSM (hunk 36)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(135, 8).Value = 99986.07000000001
$ws.Cells.Item(135, 10).Value = 99986.07000000001
$ws.Cells.Item(135, 12).Value = 99986.07000000001
$ws.Cells.Item(135, 14).Value = -110126.07

# Row 140 on GSM (hunk 37)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(140, 8).Value = 89547.37
$ws.Cells.Item(140, 10).Value = 89547.37
$ws.Cells.Item(140, 12).Value = 89547.37
$ws.Cells.Item(140, 14).Value = -99907.37

# Row 22 on LTW (hunk 38)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 52002.25
$ws.Cells.Item(22, 9).Value = 2524.5715
$ws.Cells.Item(22, 10).Value = 167450.17
$ws.Cells.Item(22, 11).Value = 2524.5715
$ws.Cells.Item(22, 12).Value = 167450.17
$ws.Cells.Item(22, 13).Value = -2229.5715
$ws.Cells.Item(22, 14).Value = -168040.17

# Row 27 on LTW (hunk 39)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 52002.25
$ws.Cells.Item(27, 9).Value = 2524.5715
$ws.Cells.Item(27, 10).Value = 167450.17
$ws.Cells.Item(27, 11).Value = 2524.5715
$ws.Cells.Item(27, 12).Value = 167450.17
$ws.Cells.Item(27, 13).Value = -2417.5715
$ws.Cells.Item(27, 14).Value = -167664.17

# Row 117 on LTW (hunk 40)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(117, 8).Value = 38988.57
$ws.Cells.Item(117, 10).Value = 38988.57
$ws.Cells.Item(117, 12).Value = 38988.57
$ws.Cells.Item(117, 14).Value = -48166.57

# Row 118 on LTW (hunk 41)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(118, 8).Value = 58811.625
$ws.Cells.Item(118, 10).Value = 58811.625
$ws.Cells.Item(118, 12).Value = 58811.625
$ws.Cells.Item(118, 14).Value = -62125.625

# Row 123 on LTW (hunk 42)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(123, 8).Value = 79996.664
$ws.Cells.Item(123, 10).Value = 79996.664
$ws.Cells.Item(123, 12).Value = 79996.664
$ws.Cells.Item(123, 14).Value = -89796.664

# Row 80 on WVR (hunk 43)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(80, 8).Value = 42550
$ws.Cells.Item(80, 10).Value = 42550
$ws.Cells.Item(80, 12).Value = 42550
$ws.Cells.Item(80, 14).Value = -44546

# Row 83 on WVR (hunk 44)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(83, 8).Value = 42550
$ws.Cells.Item(83, 10).Value = 42550
$ws.Cells.Item(83, 12).Value = 127650
$ws.Cells.Item(83, 14).Value = -137634

# Row 92 on WVR (hunk 45)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 14).ClearContents()

# Row 132 on WVR (hunk 46)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2606.611
$ws.Cells.Item(132, 9).Value = 2354.3438
$ws.Cells.Item(132, 11).Value = 7063.0314
$ws.Cells.Item(132, 13).Value = -4533.0314

# Row 136 on WVR (hunk 47)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 3275.7896
$ws.Cells.Item(136, 9).Value = 3270.6875
$ws.Cells.Item(136, 11).Value = 9812.0625
$ws.Cells.Item(136, 13).Value = -7262.0625
